$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column values that look numeric stay as text, matching the source data
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '23.969.08'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.653.59'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '308.73'
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '0.3905'
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("D8").Value = '0.3829'
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("D9").Value = '51.38'
$ws.Range("E9").Value = '  +3.53%  '
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("D11").Value = '0.9998'
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("D12").Value = '0.08444'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '24.02'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = '7.118'
$ws.Range("E14").Value = '  +0.85%  '
$ws.Range("D15").Value = '7.868'
$ws.Range("E15").Value = '  +3.88%  '
$ws.Range("D16").Value = '0.00001315'
$ws.Range("D17").Value = '1.647.63'
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").Value = '94.44'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").Value = '0.06969'
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").Value = '19.77'
$ws.Range("E20").Value = '  -1.65%  '
$ws.Range("D21").Value = '6.901'
$ws.Range("D22").Value = '0.9993'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = '13.62'
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").Value = '23.967.29'
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").Value = '2.481'
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").Value = '3.019'
$ws.Range("E26").Value = '  +5.48%  '
$ws.Range("D27").Value = '22.05'
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").Value = '153.12'
$ws.Range("E28").Value = '  -2.17%  '
$ws.Range("D29").Value = '5.434'
$ws.Range("E29").Value = '  +3.07%  '
$ws.Range("D30").Value = '139.27'
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("D31").Value = '7.733'
$ws.Range("E31").Value = '  -1.68%  '
$ws.Range("D32").Value = '2.482'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").Value = '1.830.53'
$ws.Range("E33").Value = '  +2.17%  '
$ws.Range("D34").Value = '1.036'
$ws.Range("E34").Value = '  +5.76%  '
$ws.Range("D35").Value = '0.08111'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").Value = '0.02969'
$ws.Range("E36").Value = '  +2.71%  '
$ws.Range("D37").Value = '6.754'
$ws.Range("E37").Value = '  +2.24%  '
$ws.Range("D38").Value = '10.84'
$ws.Range("E38").Value = '  +4.27%  '
$ws.Range("D39").Value = '0.2680'
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("D40").Value = '0.09157'
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").Value = '0.7557'
$ws.Range("D42").Value = '13.49'
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("D43").Value = '1.425'
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = '16.28'
$ws.Range("E44").Value = '  +1.08%  '
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("D46").Value = '2.454'
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("D48").Value = '0.9987'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").Value = '0.08304'
$ws.Range("E49").Value = '  +0.69%  '
$ws.Range("D50").Value = '134.46'
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("E51").Value = '  +1.15%  '
